$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header text: add gender-perspective note to the FAQ "Temática" column
#    header (D7).
# ---------------------------------------------------------------------------
$ws.Range("D7").Value = "Temática de las preguntas frecuentes (Redactada con perspectiva de género)"

# ---------------------------------------------------------------------------
# 2. Drop any existing hyperlinks - they will be recreated below, pointing at
#    the new target URLs (and without the old #comments anchor).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 3. Row 8 - "Obtención de título profesional electrónico"
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = 2023
$ws.Range("B8").Value = 45017
$ws.Range("C8").Value = 45107
$ws.Range("D8").Value = "Obtención de título profesional electrónico"
$ws.Range("E8").Value = "Quisiera agendar una cita para el trámite de mi título"
$ws.Range("F8").Value = @"
Hola,
Gracias por contactarnos, en el siguiente link se encuentran disponibles los Lineamientos para la Emisión y Registro de Título Profesional Electrónico:
https://www.upp.edu.mx/serviciosescolares/wp-content/uploads/2023/04/LINEAMIENTOS-ING_LTFISICA_29_03_2023.pdf
En el apartado Revisión de documentos se encuentra el procedimiento para agendar cita de revisión virtual.
Ponemos a su disposición los siguientes datos de contacto para brindarle una mejor y oportuna atención:
titulacion@upp.edu.mx
7715477510 ext. 2247
Saludos cordiales
"@
$ws.Range("G8").Value = "https://www.upp.edu.mx/serviciosescolares/?p=4156"
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = "Departamento de Servicios Escolares (UPP)"
$ws.Range("J8").Value = 45117
$ws.Range("K8").Value = 45117
$ws.Range("L8").Value = ""
$ws.Hyperlinks.Add($ws.Range("G8"), "https://www.upp.edu.mx/serviciosescolares/?p=4156")

# ---------------------------------------------------------------------------
# 4. Row 9 - "Solicitud de constancia escolar"
#    (first strip the leftover row-level custom-format flag from the old
#    content, then repaint the row's cell formats off row 8 so every cell
#    keeps the same border/alignment/wrap look as the rest of the table)
# ---------------------------------------------------------------------------
$ws.Rows.Item(9).ClearFormats()
$ws.Range("A8:L8").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A9").Value = 2023
$ws.Range("B9").Value = 45017
$ws.Range("C9").Value = 45107
$ws.Range("D9").Value = "Solicitud de constancia escolar"
$ws.Range("E9").Value = @"
Buen día para solicitar una constancia de estudios, para una persona dada de baja temporal, que ya pasaron 7 años entonces creo que ya es definitiva, me podria decir cual es el proces y si aun se puede solicitar.
es necesario matricula o solo me puedo presentar con mi nombre, es que no tengo la matricula.
Gracias, buen día
"@
$ws.Range("F9").Value = @"
Hola,
Gracias por comunicarse, en el siguiente link se encuentra disponible el procedimiento para solicitar constancia de estudios o historial académico.
https://www.upp.edu.mx/serviciosescolares/wp-content/uploads/2023/01/CATALOGO-DE-CONSTANCIAS-E-HISTORIALES-ACADEMICOS-GENERAL-30_01_2023_compressed.pdf
Ponemos a su disposición los siguientes datos de contacto para brindarle una mejor y oportuna atención:
Teléfono 7715477510 ext. 2244 y 2213
Saludos cordiales.
"@
$ws.Range("G9").Value = "https://www.upp.edu.mx/serviciosescolares/?page_id=2"
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = "Departamento de Servicios Escolares (UPP)"
$ws.Range("J9").Value = 45117
$ws.Range("K9").Value = 45117
$ws.Range("L9").Value = ""
$ws.Hyperlinks.Add($ws.Range("G9"), "https://www.upp.edu.mx/serviciosescolares/?page_id=2")

# ---------------------------------------------------------------------------
# 5. Row 10 (new) - "Información de contacto"
# ---------------------------------------------------------------------------
$ws.Range("A8:L8").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A10").Value = 2023
$ws.Range("B10").Value = 45017
$ws.Range("C10").Value = 45107
$ws.Range("D10").Value = "Información de contacto"
$ws.Range("E10").Value = @"
Hola, buenos días
Para contactarlos, favor de indicar número telefónico ya que el que está publicado en la página no está disponible…
"@
$ws.Range("F10").Value = @"
Hola,
Gracias por escribirnos, ponemos a su disposición el correo electrónico servescolares@upp.edu.mx para contactarnos.
Saludos cordiales.
"@
$ws.Range("G10").Value = "https://www.upp.edu.mx/serviciosescolares/?page_id=2"
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = "Departamento de Servicios Escolares (UPP)"
$ws.Range("J10").Value = 45117
$ws.Range("K10").Value = 45117
$ws.Range("L10").Value = ""
$ws.Hyperlinks.Add($ws.Range("G10"), "https://www.upp.edu.mx/serviciosescolares/?page_id=2")

# ---------------------------------------------------------------------------
# 6. Row 11 (new) - "Información de nuevo ingreso"
# ---------------------------------------------------------------------------
$ws.Range("A8:L8").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = 2023
$ws.Range("B11").Value = 45017
$ws.Range("C11").Value = 45107
$ws.Range("D11").Value = "Información de nuevo ingreso"
$ws.Range("E11").Value = "Hola me gustaría saber la fecha del examen de admisión en donde puedo consultarlo gracias"
$ws.Range("F11").Value = @"
Hola,
Gracias por escribirnos, se pide atentamente envíe un correo a servescolares@upp.edu.mx con los siguientes datos:
Nombre completo
CURP
Número de ficha
Programa educativo (carrera) a la que desea ingresar
Saludos cordiales
"@
$ws.Range("G11").Value = "https://www.upp.edu.mx/serviciosescolares/?p=4674"
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = "Departamento de Servicios Escolares (UPP)"
$ws.Range("J11").Value = 45117
$ws.Range("K11").Value = 45117
$ws.Range("L11").Value = ""
$ws.Hyperlinks.Add($ws.Range("G11"), "https://www.upp.edu.mx/serviciosescolares/?p=4674")

# ---------------------------------------------------------------------------
# 7. Row heights (autofit-by-hand, matching the taller wrapped text) and the
#    header row's slightly shorter height.
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 38.25
$ws.Rows.Item(8).RowHeight = 180
$ws.Rows.Item(9).RowHeight = 150
$ws.Rows.Item(10).RowHeight = 88.5
$ws.Rows.Item(11).RowHeight = 161.25

# ---------------------------------------------------------------------------
# 8. Column widths - widened to fit the new, longer answer text.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 65.31
$ws.Columns.Item(5).ColumnWidth = 60.17
$ws.Columns.Item(6).ColumnWidth = 73.17
$ws.Columns.Item(7).ColumnWidth = 62.17

# ---------------------------------------------------------------------------
# 9. View: page orientation + active selection.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("B8").Select()
